$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.818.12"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +4.42%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.844.22"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +6.50%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.58%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'423.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +4.95%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'129.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.25%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.834.72"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +6.50%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -1.18%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -0.23%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.722"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +0.33%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.157"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.74%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.0000334"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +5.53%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'40.82"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.35%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'10.42"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +5.89%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'4.456.38"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +5.82%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'15.71"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +17.04%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.856.60"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +6.66%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'  -0.80%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'19.88"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.43%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'66.957.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +3.93%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  +1.00%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'413.69"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.95%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'14.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.39%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'84.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.95%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +2.15%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'37.58"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +6.31%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +7.02%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'3.25"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +2.31%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'9.32"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +35.16%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'5.36"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +4.48%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'735.52"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +8.94%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'13.06"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +2.75%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.122"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +5.56%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +2.40%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.998"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.04%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -5.15%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'38.64"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -5.52%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'55.58"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.03%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'5.48"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +24.94%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("B40").Value = "'PEPE"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = "'https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'0.0₃0727"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +14.74%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = "'VeChain"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'0.0457"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.78%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'2.89"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.21%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +0.49%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +2.49%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -3.66%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.318"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +10.22%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -0.65%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'2.05"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.48%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'140.79"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -2.49%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'2.81"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.58%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +0.36%  "
$ws.Range("E51").Style = "Normal"
